# Turn the bare "https://www.youtube.com/watch?v=r9buAwVBDhA" paragraph into
# a real hyperlink (styled with the "Hyperlink" character style), then add a
# new paragraph right after it containing a second hyperlink to the
# khinsider album page, surrounded by the same blank-paragraph spacing shown
# in the target revision.

$d = $word.ActiveDocument

$youtubeUrl   = "https://www.youtube.com/watch?v=r9buAwVBDhA"
$khinsiderUrl = "https://downloads.khinsider.com/game-soundtracks/album/happy-mario-20th-super-mario-sound-collection"

# --- Locate the paragraph that still holds the plain-text YouTube link -----
$ytIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.TrimEnd() -eq $youtubeUrl) {
        $ytIndex = $i
        break
    }
}
if ($ytIndex -eq -1) {
    throw "Could not find the YouTube URL paragraph"
}

# --- Step 1: wrap that run in a hyperlink field -----------------------------
$ytPara  = $d.Paragraphs($ytIndex)
$ytRange = $ytPara.Range
$ytRange.End = $ytRange.End - 1   # exclude the paragraph mark
$d.Hyperlinks.Add($ytRange, $youtubeUrl) | Out-Null

# --- Step 2: insert a fresh, empty paragraph right after it ----------------
$ytPara     = $d.Paragraphs($ytIndex)
$insertion  = $ytPara.Range
$insertion.Collapse(0)
$insertion.InsertParagraphAfter()

# The freshly inserted paragraph now sits at $ytIndex + 1 and contains a
# single empty run. Stamp it with unique placeholder text so later Find
# operations can target exactly that paragraph, then grow the surrounding
# blank paragraphs one `Find.Execute` pass at a time (each pass only ever
# inserts a single "^p" so Word collapses the run back down to a clean
# empty <w:p/> instead of leaving a stray empty run behind).
$newIndex = $ytIndex + 1
$placeholder = "KHINSIDERPLACEHOLDER"
$newPara = $d.Paragraphs($newIndex)
$newRange = $newPara.Range
$newRange.End = $newRange.End - 1
$newRange.Text = $placeholder

# Leading blank paragraph before the khinsider text.
$d.Content.Find.Execute($placeholder, $true, $false, $false, $false, $false, `
    $true, 1, $false, "^p$placeholder", 2) | Out-Null

# Two trailing blank paragraphs after the khinsider text (done as two
# separate single-"^p" passes — see note above).
$d.Content.Find.Execute($placeholder, $true, $false, $false, $false, $false, `
    $true, 1, $false, "$placeholder^p", 2) | Out-Null
$d.Content.Find.Execute($placeholder, $true, $false, $false, $false, $false, `
    $true, 1, $false, "$placeholder^p", 2) | Out-Null

# Swap the placeholder for the real URL text.
$d.Content.Find.Execute($placeholder, $true, $false, $false, $false, $false, `
    $true, 1, $false, $khinsiderUrl, 2) | Out-Null

# --- Step 3: hyperlink the khinsider paragraph ------------------------------
$khIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.TrimEnd() -eq $khinsiderUrl) {
        $khIndex = $i
        break
    }
}
if ($khIndex -eq -1) {
    throw "Could not find the khinsider URL paragraph"
}

$khPara  = $d.Paragraphs($khIndex)
$khRange = $khPara.Range
$khRange.End = $khRange.End - 1   # exclude the paragraph mark
$d.Hyperlinks.Add($khRange, $khinsiderUrl) | Out-Null
